$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '43.607.33'
$ws.Range('D3').Value = '2.290.90'
$ws.Range('E3').Value = '  +0.79%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '96.44'
$ws.Range('E5').Value = '  +2.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '266.67'
$ws.Range('E6').Value = '  -0.96%  '
$ws.Range('E7').Value = '  -0.43%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.609'
$ws.Range('E9').Value = '  -2.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '45.69'
$ws.Range('E10').Value = '  -1.97%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0935'
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('E12').Value = '  -2.89%  '
$ws.Range('E13').Value = '  +0.26%  '
$ws.Range('D14').Value = '2.637.09'
$ws.Range('E14').Value = '  +0.90%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.16'
$ws.Range('E15').Value = '  -1.29%  '
$ws.Range('E16').Value = '  +2.48%  '
$ws.Range('D17').Value = '2.295.60'
$ws.Range('E17').Value = '  +0.46%  '
$ws.Range('D18').Value = '43.629.72'
$ws.Range('E18').Value = '  -0.75%  '
$ws.Range('E19').Value = '  +2.00%  '
$ws.Range('E20').Value = '  +0.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.72'
$ws.Range('E21').Value = '  +1.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.43'
$ws.Range('E22').Value = '  +5.97%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '232.18'
$ws.Range('E23').Value = '  -1.63%  '
$ws.Range('E24').Value = '  -9.28%  '
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.15'
$ws.Range('E27').Value = '  -1.59%  '
$ws.Range('E28').Value = '  +3.38%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '40.10'
$ws.Range('E29').Value = '  +1.97%  '
$ws.Range('E30').Value = '  +0.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '175.46'
$ws.Range('E31').Value = '  +1.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.83'
$ws.Range('E32').Value = '  -1.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0880'
$ws.Range('E33').Value = '  -3.00%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.36'
$ws.Range('E34').Value = '  -4.06%  '
$ws.Range('E35').Value = '  +0.76%  '
$ws.Range('E36').Value = '  -2.71%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0354'
$ws.Range('E37').Value = '  +1.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.33'
$ws.Range('E38').Value = '  -2.15%  '
$ws.Range('E39').Value = '  -2.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.235'
$ws.Range('E40').Value = '  -4.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.31'
$ws.Range('E41').Value = '  +3.81%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '12.23'
$ws.Range('E42').Value = '  -0.87%  '
$ws.Range('E43').Value = '  +10.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '64.16'
$ws.Range('E44').Value = '  +4.42%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.80'
$ws.Range('E45').Value = '  +3.17%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.22'
$ws.Range('E46').Value = '  -4.54%  '
$ws.Range('E47').Value = '  -0.22%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '97.93'
$ws.Range('E48').Value = '  -1.43%  '
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('D50').Value = '2.515.07'
$ws.Range('E50').Value = '  +0.79%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.426'
$ws.Range('E51').Value = '  -0.59%  '
